# edit.ps1 — applies the commit's changes:
#   1. Re-style the three tables (slides 14, 15, 16) from table style
#      {4D87F191-96BA-49B0-B6BF-1D5145A105D7} to
#      {90D82F63-57D0-4330-ADD3-6C407542997A}.
#   2. Swap the deck's theme palette: the slide master's theme ("Integral" /
#      "Red Violet") becomes the stock Office theme colours ("Office").

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------------
$newTableStyleId = "{90D82F63-57D0-4330-ADD3-6C407542997A}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colour scheme -------------------------------------------------
$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
